$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-06-10 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-06-11 Sunday", 2) | Out-Null
$d.Content.Find.Execute("3+5=", $true, $false, $false, $false, $false, $true, 1, $false, "42-28=", 2) | Out-Null
$d.Content.Find.Execute("80-39=", $true, $false, $false, $false, $false, $true, 1, $false, "42-33=", 2) | Out-Null
$d.Content.Find.Execute("58-12=", $true, $false, $false, $false, $false, $true, 1, $false, "90-52=", 2) | Out-Null
$d.Content.Find.Execute("10+39=", $true, $false, $false, $false, $false, $true, 1, $false, "98-24=", 2) | Out-Null
$d.Content.Find.Execute("14+41=", $true, $false, $false, $false, $false, $true, 1, $false, "1+39=", 2) | Out-Null
$d.Content.Find.Execute("95-36=", $true, $false, $false, $false, $false, $true, 1, $false, "94-57=", 2) | Out-Null
$d.Content.Find.Execute("47-29=", $true, $false, $false, $false, $false, $true, 1, $false, "76-22=", 2) | Out-Null
$d.Content.Find.Execute("96-16=", $true, $false, $false, $false, $false, $true, 1, $false, "52-3=", 2) | Out-Null
$d.Content.Find.Execute("55-20=", $true, $false, $false, $false, $false, $true, 1, $false, "37-6=", 2) | Out-Null
$d.Content.Find.Execute("30+41=", $true, $false, $false, $false, $false, $true, 1, $false, "8+53=", 2) | Out-Null
$d.Content.Find.Execute("18+56=", $true, $false, $false, $false, $false, $true, 1, $false, "94-65=", 2) | Out-Null
$d.Content.Find.Execute("47-13=", $true, $false, $false, $false, $false, $true, 1, $false, "75-39=", 2) | Out-Null
$d.Content.Find.Execute("26+26=", $true, $false, $false, $false, $false, $true, 1, $false, "94-21=", 2) | Out-Null
$d.Content.Find.Execute("95-43=", $true, $false, $false, $false, $false, $true, 1, $false, "44+2=", 2) | Out-Null
$d.Content.Find.Execute("9+12=", $true, $false, $false, $false, $false, $true, 1, $false, "28+57=", 2) | Out-Null
$d.Content.Find.Execute("16-7=", $true, $false, $false, $false, $false, $true, 1, $false, "99-77=", 2) | Out-Null
$d.Content.Find.Execute("91-29=", $true, $false, $false, $false, $false, $true, 1, $false, "70+2=", 2) | Out-Null
$d.Content.Find.Execute("34-19=", $true, $false, $false, $false, $false, $true, 1, $false, "80-22=", 2) | Out-Null
$d.Content.Find.Execute("43-31=", $true, $false, $false, $false, $false, $true, 1, $false, "95-0=", 2) | Out-Null
$d.Content.Find.Execute("69-53=", $true, $false, $false, $false, $false, $true, 1, $false, "41-7=", 2) | Out-Null
$d.Content.Find.Execute("50+37=", $true, $false, $false, $false, $false, $true, 1, $false, "20+35=", 2) | Out-Null
$d.Content.Find.Execute("6-3=", $true, $false, $false, $false, $false, $true, 1, $false, "70-42=", 2) | Out-Null
$d.Content.Find.Execute("69+15=", $true, $false, $false, $false, $false, $true, 1, $false, "74-48=", 2) | Out-Null
$d.Content.Find.Execute("62-53=", $true, $false, $false, $false, $false, $true, 1, $false, "72-3=", 2) | Out-Null
$d.Content.Find.Execute("56-51=", $true, $false, $false, $false, $false, $true, 1, $false, "64-7=", 2) | Out-Null
$d.Content.Find.Execute("45+1=", $true, $false, $false, $false, $false, $true, 1, $false, "21+22=", 2) | Out-Null
$d.Content.Find.Execute("96-49=", $true, $false, $false, $false, $false, $true, 1, $false, "7+29=", 2) | Out-Null
$d.Content.Find.Execute("49-48=", $true, $false, $false, $false, $false, $true, 1, $false, "32-5=", 2) | Out-Null
$d.Content.Find.Execute("82-23=", $true, $false, $false, $false, $false, $true, 1, $false, "38+3=", 2) | Out-Null
$d.Content.Find.Execute("76-48=", $true, $false, $false, $false, $false, $true, 1, $false, "97-68=", 2) | Out-Null
$d.Content.Find.Execute("29+26=", $true, $false, $false, $false, $false, $true, 1, $false, "95-19=", 2) | Out-Null
$d.Content.Find.Execute("91-38=", $true, $false, $false, $false, $false, $true, 1, $false, "53+14=", 2) | Out-Null
$d.Content.Find.Execute("90-7=", $true, $false, $false, $false, $false, $true, 1, $false, "67-31=", 2) | Out-Null
$d.Content.Find.Execute("61+26=", $true, $false, $false, $false, $false, $true, 1, $false, "64-0=", 2) | Out-Null
$d.Content.Find.Execute("34+0=", $true, $false, $false, $false, $false, $true, 1, $false, "0+56=", 2) | Out-Null
$d.Content.Find.Execute("12+85=", $true, $false, $false, $false, $false, $true, 1, $false, "1+30=", 2) | Out-Null
$d.Content.Find.Execute("26+48=", $true, $false, $false, $false, $false, $true, 1, $false, "44+51=", 2) | Out-Null
$d.Content.Find.Execute("22-7=", $true, $false, $false, $false, $false, $true, 1, $false, "7+14=", 2) | Out-Null
$d.Content.Find.Execute("3+61=", $true, $false, $false, $false, $false, $true, 1, $false, "93+2=", 2) | Out-Null
$d.Content.Find.Execute("88-69=", $true, $false, $false, $false, $false, $true, 1, $false, "86-75=", 2) | Out-Null
$d.Content.Find.Execute("19+41=", $true, $false, $false, $false, $false, $true, 1, $false, "7+24=", 2) | Out-Null
$d.Content.Find.Execute("42+43=", $true, $false, $false, $false, $false, $true, 1, $false, "59-19=", 2) | Out-Null
$d.Content.Find.Execute("91-22=", $true, $false, $false, $false, $false, $true, 1, $false, "51-20=", 2) | Out-Null
$d.Content.Find.Execute("72-24=", $true, $false, $false, $false, $false, $true, 1, $false, "96-68=", 2) | Out-Null
$d.Content.Find.Execute("13-4=", $true, $false, $false, $false, $false, $true, 1, $false, "42-30=", 2) | Out-Null
$d.Content.Find.Execute("80+17=", $true, $false, $false, $false, $false, $true, 1, $false, "10+35=", 2) | Out-Null
$d.Content.Find.Execute("49+24=", $true, $false, $false, $false, $false, $true, 1, $false, "28-19=", 2) | Out-Null
$d.Content.Find.Execute("95-31=", $true, $false, $false, $false, $false, $true, 1, $false, "36-24=", 2) | Out-Null
$d.Content.Find.Execute("65+25=", $true, $false, $false, $false, $false, $true, 1, $false, "99-22=", 2) | Out-Null
$d.Content.Find.Execute("87+6=", $true, $false, $false, $false, $false, $true, 1, $false, "9+39=", 2) | Out-Null
$d.Content.Find.Execute("71-21=", $true, $false, $false, $false, $false, $true, 1, $false, "50-14=", 2) | Out-Null
$d.Content.Find.Execute("73+15=", $true, $false, $false, $false, $false, $true, 1, $false, "4+61=", 2) | Out-Null
$d.Content.Find.Execute("60-30=", $true, $false, $false, $false, $false, $true, 1, $false, "75-59=", 2) | Out-Null
$d.Content.Find.Execute("40+14=", $true, $false, $false, $false, $false, $true, 1, $false, "66-46=", 2) | Out-Null
$d.Content.Find.Execute("8+5=", $true, $false, $false, $false, $false, $true, 1, $false, "46+35=", 2) | Out-Null
$d.Content.Find.Execute("44+3=", $true, $false, $false, $false, $false, $true, 1, $false, "75-42=", 2) | Out-Null
$d.Content.Find.Execute("9+70=", $true, $false, $false, $false, $false, $true, 1, $false, "22+3=", 2) | Out-Null
$d.Content.Find.Execute("63-29=", $true, $false, $false, $false, $false, $true, 1, $false, "46+44=", 2) | Out-Null
$d.Content.Find.Execute("22+6=", $true, $false, $false, $false, $false, $true, 1, $false, "96-39=", 2) | Out-Null
$d.Content.Find.Execute("93+1=", $true, $false, $false, $false, $false, $true, 1, $false, "75-19=", 2) | Out-Null
$d.Content.Find.Execute("39+38=", $true, $false, $false, $false, $false, $true, 1, $false, "27-0=", 2) | Out-Null
$d.Content.Find.Execute("46-42=", $true, $false, $false, $false, $false, $true, 1, $false, "69-66=", 2) | Out-Null
$d.Content.Find.Execute("84-2=", $true, $false, $false, $false, $false, $true, 1, $false, "55-50=", 2) | Out-Null
$d.Content.Find.Execute("28-10=", $true, $false, $false, $false, $false, $true, 1, $false, "76-71=", 2) | Out-Null
$d.Content.Find.Execute("36+39=", $true, $false, $false, $false, $false, $true, 1, $false, "25-24=", 2) | Out-Null
$d.Content.Find.Execute("82-24=", $true, $false, $false, $false, $false, $true, 1, $false, "15+22=", 2) | Out-Null
$d.Content.Find.Execute("79-43=", $true, $false, $false, $false, $false, $true, 1, $false, "46+15=", 2) | Out-Null
$d.Content.Find.Execute("75-51=", $true, $false, $false, $false, $false, $true, 1, $false, "93-64=", 2) | Out-Null
$d.Content.Find.Execute("91+4=", $true, $false, $false, $false, $false, $true, 1, $false, "65-43=", 2) | Out-Null
$d.Content.Find.Execute("90-48=", $true, $false, $false, $false, $false, $true, 1, $false, "41+6=", 2) | Out-Null
$d.Content.Find.Execute("59-38=", $true, $false, $false, $false, $false, $true, 1, $false, "78-28=", 2) | Out-Null
$d.Content.Find.Execute("34+52=", $true, $false, $false, $false, $false, $true, 1, $false, "29-25=", 2) | Out-Null
$d.Content.Find.Execute("9+25=", $true, $false, $false, $false, $false, $true, 1, $false, "77-71=", 2) | Out-Null
$d.Content.Find.Execute("84-42=", $true, $false, $false, $false, $false, $true, 1, $false, "98-52=", 2) | Out-Null
$d.Content.Find.Execute("25+48=", $true, $false, $false, $false, $false, $true, 1, $false, "74-48=", 2) | Out-Null
$d.Content.Find.Execute("16+13=", $true, $false, $false, $false, $false, $true, 1, $false, "33-2=", 2) | Out-Null
$d.Content.Find.Execute("27+1=", $true, $false, $false, $false, $false, $true, 1, $false, "70-53=", 2) | Out-Null
$d.Content.Find.Execute("61-19=", $true, $false, $false, $false, $false, $true, 1, $false, "86-22=", 2) | Out-Null
$d.Content.Find.Execute("76-62=", $true, $false, $false, $false, $false, $true, 1, $false, "90-6=", 2) | Out-Null
$d.Content.Find.Execute("98+0=", $true, $false, $false, $false, $false, $true, 1, $false, "83-23=", 2) | Out-Null
$d.Content.Find.Execute("49+49=", $true, $false, $false, $false, $false, $true, 1, $false, "89-87=", 2) | Out-Null
$d.Content.Find.Execute("55-23=", $true, $false, $false, $false, $false, $true, 1, $false, "7+89=", 2) | Out-Null
$d.Content.Find.Execute("16+60=", $true, $false, $false, $false, $false, $true, 1, $false, "87-55=", 2) | Out-Null
$d.Content.Find.Execute("26+34=", $true, $false, $false, $false, $false, $true, 1, $false, "49-4=", 2) | Out-Null
$d.Content.Find.Execute("95-85=", $true, $false, $false, $false, $false, $true, 1, $false, "80+14=", 2) | Out-Null
$d.Content.Find.Execute("45+28=", $true, $false, $false, $false, $false, $true, 1, $false, "93-79=", 2) | Out-Null
$d.Content.Find.Execute("70-28=", $true, $false, $false, $false, $false, $true, 1, $false, "12+43=", 2) | Out-Null
$d.Content.Find.Execute("48-26=", $true, $false, $false, $false, $false, $true, 1, $false, "35+28=", 2) | Out-Null
$d.Content.Find.Execute("84-40=", $true, $false, $false, $false, $false, $true, 1, $false, "38+1=", 2) | Out-Null
$d.Content.Find.Execute("21+59=", $true, $false, $false, $false, $false, $true, 1, $false, "85-74=", 2) | Out-Null
$d.Content.Find.Execute("12+71=", $true, $false, $false, $false, $false, $true, 1, $false, "99-65=", 2) | Out-Null
$d.Content.Find.Execute("17+53=", $true, $false, $false, $false, $false, $true, 1, $false, "35+35=", 2) | Out-Null
$d.Content.Find.Execute("99-13=", $true, $false, $false, $false, $false, $true, 1, $false, "71-17=", 2) | Out-Null
$d.Content.Find.Execute("36+45=", $true, $false, $false, $false, $false, $true, 1, $false, "18-16=", 2) | Out-Null
$d.Content.Find.Execute("80+7=", $true, $false, $false, $false, $false, $true, 1, $false, "48+29=", 2) | Out-Null
$d.Content.Find.Execute("61-51=", $true, $false, $false, $false, $false, $true, 1, $false, "53+24=", 2) | Out-Null
$d.Content.Find.Execute("61+19=", $true, $false, $false, $false, $false, $true, 1, $false, "77-11=", 2) | Out-Null
$d.Content.Find.Execute("18-14=", $true, $false, $false, $false, $false, $true, 1, $false, "4+87=", 2) | Out-Null
$d.Content.Find.Execute("15+64=", $true, $false, $false, $false, $false, $true, 1, $false, "31+25=", 2) | Out-Null
$d.Content.Find.Execute("24+20=", $true, $false, $false, $false, $false, $true, 1, $false, "43+33=", 2) | Out-Null
